$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

# Column A: index number - copy formatting (bold/border/centered) from A96, then set value
$ws.Cells.Item(96, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = 96

$ws.Cells.Item($row, 2).Value = "denmark"
$ws.Cells.Item($row, 3).Value = "superliga"
$ws.Cells.Item($row, 4).Value = "2023-2024"

# Column E: match date - copy formatting (custom date number format) from E96, then set value
$ws.Cells.Item(96, 5).Copy($ws.Cells.Item($row, 5))
$ws.Cells.Item($row, 5).Value = 45257.79166666666

$ws.Cells.Item($row, 6).Value = "Silkeborg"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Midtjylland"
$ws.Cells.Item($row, 9).Value = 4
$ws.Cells.Item($row, 10).Value = 2.44
$ws.Cells.Item($row, 11).Value = "13/11/2023 10:42"
$ws.Cells.Item($row, 12).Value = 3.16
$ws.Cells.Item($row, 13).Value = "27/11/2023 18:59"
$ws.Cells.Item($row, 14).Value = 3.35
$ws.Cells.Item($row, 15).Value = "13/11/2023 10:42"
$ws.Cells.Item($row, 16).Value = 3.48
$ws.Cells.Item($row, 17).Value = "27/11/2023 18:59"
$ws.Cells.Item($row, 18).Value = 2.95
$ws.Cells.Item($row, 19).Value = "13/11/2023 10:42"
$ws.Cells.Item($row, 20).Value = 2.31
$ws.Cells.Item($row, 21).Value = "27/11/2023 17:54"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/silkeborg-midtjylland/tCaJGNSk/"
